$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay as plain text, matching the
# original inline-string cell contents (prices use dotted "thousand" groupings
# and leading zeros that Excel would otherwise mangle by auto-converting to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.191.20"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.802.26"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "336.79"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "0.9972"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "0.3958"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("D8").Value = "0.3442"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "45.42"
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("D10").Value = "1.151"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").Value = "0.07380"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "22.87"
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("D13").Value = "0.9964"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "6.248"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "7.289"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "1.796.46"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "0.00001081"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "0.06634"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "81.42"
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("D20").Value = "0.9970"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").Value = "6.318"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "28.187.37"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "11.86"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").Value = "2.396"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "20.49"
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.419"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "155.29"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "1.993.40"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "1.323"
$ws.Range("E30").Value = "  -7.33%  "
$ws.Range("D31").Value = "131.07"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").Value = "4.059"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "5.970"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").Value = "0.08753"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("D35").Value = "12.34"
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("D36").Value = "0.06315"
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").Value = "0.6669"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02325"
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.228"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").Value = "0.2143"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").Value = "1.513"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").Value = "1.214"
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("D43").Value = "8.139"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "14.16"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "0.9971"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "0.6138"
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").Value = "3.847"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "128.38"
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").Value = "2.036"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").Value = "1.172"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").Value = "0.07089"
$ws.Range("E51").Value = "  -5.29%  "
